# reworked-fig7.pptx -- "some work on Intro and grouping" commit.
#
# This deck's single slide holds the big multi-panel "replace()" figure.
# The edit:
#   1) Nudges the small "Group 330" group (the oval/rectangle cluster near
#      the bottom-left of the figure) a little left/down.
#   2) Straightens out the arrow connector feeding into that group: it is
#      no longer vertically flipped, and its off/ext change slightly to
#      match the new endpoints.
#   3) Adds 15 new small caption textboxes ("wgBy" x3, "used" x12) that
#      label pieces of the figure -- these are appended at the end of the
#      shape tree, after the trailing "replace()" textbox.
#
# (NB: the same commit also bumped a bunch of footer date fields
# 9/9/19 -> 9/27/19 on other slides of the original multi-slide deck, but
# this working copy only has the one slide and it carries no date-field
# placeholders, so there is nothing to do for that part here.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $null
}

$EMU_PER_POINT = 12700.0

# ---------------------------------------------------------------------
# 1) Reposition the "Group 330" group.
# ---------------------------------------------------------------------
$grp = Find-ShapeByName $s.Shapes "Group 330"
$grp.Left = 466872 / $EMU_PER_POINT
$grp.Top = 3751439 / $EMU_PER_POINT

# ---------------------------------------------------------------------
# 2) Fix up the "Straight Arrow Connector 335" feeding into it.
# ---------------------------------------------------------------------
$cxn = Find-ShapeByName $s.Shapes "Straight Arrow Connector 335"
$cxn.VerticalFlip = $false
$cxn.Left = 837330 / $EMU_PER_POINT
$cxn.Top = 3935742 / $EMU_PER_POINT
$cxn.Width = 1525315 / $EMU_PER_POINT
$cxn.Height = 12963 / $EMU_PER_POINT

# ---------------------------------------------------------------------
# 3) Add the new caption textboxes.
# ---------------------------------------------------------------------
$newLabels = @(
    @{ Name = "TextBox 273"; X = 847593;  Y = 1096378; W = 510333; H = 276999; Text = "wgBy" },
    @{ Name = "TextBox 276"; X = 1794365; Y = 2245440; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 280"; X = 1895314; Y = 1380858; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 282"; X = 1747182; Y = 580338;  W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 283"; X = 7353139; Y = 2222175; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 285"; X = 7408513; Y = 1386041; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 286"; X = 7307361; Y = 646469;  W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 287"; X = 3518887; Y = 1249328; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 288"; X = 8979011; Y = 1275604; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 289"; X = 8967964; Y = 3772236; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 290"; X = 7231326; Y = 3140651; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 291"; X = 7207512; Y = 3794317; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 292"; X = 7314635; Y = 4684269; W = 482824; H = 276999; Text = "used" },
    @{ Name = "TextBox 293"; X = 6291889; Y = 3607574; W = 510333; H = 276999; Text = "wgBy" },
    @{ Name = "TextBox 294"; X = 6339465; Y = 1118760; W = 510333; H = 276999; Text = "wgBy" }
)

foreach ($label in $newLabels) {
    $tb = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $tb.Name = $label.Name

    $tf = $tb.TextFrame
    $tf.WordWrap = $false
    $tf.AutoSize = 1

    $tr = $tf.TextRange
    $tr.Text = $label.Text
    $tr.Font.Size = 12
    $tr.Font.Name = "+mj-lt"
    $tb.TextFrame2.TextRange.Font.NameComplexScript = "Times New Roman"

    $tb.Fill.Visible = $false

    $tb.Left = $label.X / $EMU_PER_POINT
    $tb.Top = $label.Y / $EMU_PER_POINT
    $tb.Width = $label.W / $EMU_PER_POINT
    $tb.Height = $label.H / $EMU_PER_POINT
}
